$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Row 15: battery spec capacity changed from 370Wh to 180Wh
$ws.Range("B15").Value = "180Wh"

# Row 19: "Estimated capacity" now derives from the new 180Wh figure;
# keep the formula (=180/0.9) in B19 and mark the empty D19 cell with the
# same percentage format as the other "depth of discharge" style cell.
$ws.Range("B19").Formula = "=180/0.9"
$ws.Range("D19").NumberFormat = "0%"

# Row 20 ("Size (litre)"): B20 becomes a computed estimate, the old measured
# value moves to C20, and D20 is labelled "est" to flag it as an estimate.
$ws.Range("B20").Formula = "=2*90*46*158*10^-9 * 10^3"
$ws.Range("C20").Value = 1.1200000000000001
$ws.Range("D20").Value = "est"

# Row 21 ("Mass (kg)"): same pattern - B21 becomes a computed estimate, the
# old measured value moves to C21, D21 labelled "est".
$ws.Range("B21").Formula = "=1242*2*10^-3"
$ws.Range("C21").Value = 2.7
$ws.Range("D21").Value = "est"

# Column B now holds mixed numeric/text content across rows 15-21; autofit
# it like the new layout does.
$ws.Columns.Item(2).AutoFit()

# Selection moves to the new last row used (B21).
$ws.Range("B21").Select()
